# "Priyanka - pdf formate"
# Update the October-2014 bank statement sheet:
#   - Row 2 (Sekhar Beri): fill in the account number and correct Netpay
#   - Row 3 (BalaRaju Vankala): fill in an account-number placeholder and correct Netpay
#   - Row 4 (Priyanka Muddana): correct Netpay
#   - Row 5 (pattabhi ramarao) is removed entirely
#   - Column A is widened to match column B

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Sekhar Beri
$ws.Range("A2").Value = 17249172304
$ws.Range("C2").Value = 1304.86

# Row 3 - BalaRaju Vankala
$ws.Range("A3").Value = "awetwetawe"
$ws.Range("C3").Value = 1449.85

# Row 4 - Priyanka Muddana (Netpay only; account/name/month untouched)
$ws.Range("C4").Value = 17057.4

# Row 5 (pattabhi ramarao) no longer exists in the updated statement
$ws.Rows.Item(5).Delete()

# Column A widens to line up with column B
$ws.Columns.Item(1).ColumnWidth = 14.43
